# Update "想去人数" (want-to-go count) figures in F column for both the
# "展览" and "全部类型" worksheets, reflecting newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2-6
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3243
$wsExhibit.Range("F3").Value = 2
$wsExhibit.Range("F4").Value = 53
$wsExhibit.Range("F5").Value = 1143
$wsExhibit.Range("F6").Value = 307

# Sheet "全部类型" - rows 2-5 and 7 (row 6 unchanged)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3243
$wsAll.Range("F3").Value = 2
$wsAll.Range("F4").Value = 53
$wsAll.Range("F5").Value = 1143
$wsAll.Range("F7").Value = 307
